$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Fonctionnalité"
$ws.Range("B1").Value = "Description"
$ws.Range("A3").Value = "flux de questions"
$ws.Range("B3").Value = "affiche l'ensemble des questions"
$ws.Range("A4").Value = "ajouter question"
$ws.Range("B4").Value = "ajouter une question stockée en base de données"
$ws.Range("A5").Value = "supprimer question"
$ws.Range("B5").Value = "supprimer une question stockée en base de données"
$ws.Range("A6").Value = "modification question"
$ws.Range("B6").Value = "modifier une question stockée en base de données"
$ws.Range("A7").Value = "rechercher question"
$ws.Range("B7").Value = "rechercher une question par son titre"
$ws.Range("A8").Value = "répondre question"
$ws.Range("B8").Value = "répondre à une question "
$ws.Range("A9").Value = "supprimer réponse"
$ws.Range("B9").Value = "supprimer sa réponse à une question"
$ws.Range("A10").Value = "question privé"
$ws.Range("B10").Value = "question uniquement visible par ses amis"
$ws.Range("A11").Value = "aimer question"
$ws.Range("B11").Value = "aimer une question"
$ws.Range("A12").Value = "je n'aime pas question"
$ws.Range("B12").Value = "ne pas aimer une question"
$ws.Range("A13").Value = "retirer aime question"
$ws.Range("B13").Value = "retirer le j'aime d'une question"
$ws.Range("A14").Value = "modification profil"
$ws.Range("B14").Value = "modifier ses données personnel"
$ws.Range("A15").Value = "inscription utilisateur"
$ws.Range("B15").Value = "formulaire pour s'incrire"
$ws.Range("A16").Value = "connexion utilisateur"
$ws.Range("B16").Value = "formulaire pour se connecter"
$ws.Range("A17").Value = "vue question"
$ws.Range("B17").Value = "voir une question en fonction de son id"
$ws.Range("A18").Value = "vue utilisateur "
$ws.Range("B18").Value = "voir un utilisateur en fonction de son id"
$ws.Range("A19").Value = "ajouter ami"
$ws.Range("B19").Value = "ajouter un utilisateur en ami"
$ws.Range("A20").Value = "supprimer ami"
$ws.Range("B20").Value = "supprimer un utilisateur en ami"
$ws.Range("A21").Value = "voir ami"
$ws.Range("B21").Value = "voir la liste de ses amis"
$ws.Range("A22").Value = "espace administrateur"
$ws.Range("B22").Value = "back end pour gérer le site"
$ws.Range("A23").Value = "supprimer utilisateur"
$ws.Range("B23").Value = "permet de supprimer définitivement un utilisateur quand t-on est admin"
$ws.Range("A24").Value = "supprimer son compte"
$ws.Range("B24").Value = "permet de supprimer son compte"

$ws.Columns("A:A").AutoFit() | Out-Null
$ws.Range("C1").Select() | Out-Null
